$d = $word.ActiveDocument

# Step 1: Replace paragraph 1 (strip the _GoBack bookmark)
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:u w:val=`"single`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:u w:val=`"single`"/></w:rPr><w:t>Publication-Related Comments</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/><w:u w:val=`"single`"/></w:rPr><w:t xml:space=`"preserve`">: </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Step 2: Replace paragraph 2 (merge runs: "Graham: " + "To discuss..." ; merge "...merits " + "authorship." dropping its proofErr wrap)
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">Graham: To discuss later…..I believe </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>Lianna</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> will want to be an author and, if that is justified then so is making Linda an author.  If we do this, it well be if they contribute to future drafts of the MS, and if we all agree after a discussion among the 4 of that everyone’s contribution merits authorship.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Step 3: Append two new list paragraphs at the end of the document body (before sectPr),
# moving the _GoBack bookmark into the new final paragraph.
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
$endRange.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>Graham: RE fig. 8 (true fish richness model)</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">`"Not important now, but I would later be curious to see this graph with the points </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>colour</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> coded by time`"</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
